$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.324.98"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "3.138.24"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'609.89"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").Value = "'143.52"
$ws.Range("E6").Value = "  -2.57%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.134.63"
$ws.Range("E8").Value = "  -0.57%  "
$ws.Range("D9").Value = "'0.529"
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").Value = "'5.36"
$ws.Range("E11").Value = "  -3.02%  "
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("E13").Value = "  +2.64%  "
$ws.Range("D14").Value = "'35.50"
$ws.Range("E14").Value = "  -1.71%  "
$ws.Range("D15").Value = "3.650.70"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("E16").Value = "  +2.58%  "
$ws.Range("D17").Value = "64.292.60"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").Value = "3.155.85"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("E19").Value = "  -1.31%  "
$ws.Range("D20").Value = "'477.99"
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("D21").Value = "'14.70"
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("D22").Value = "'0.721"
$ws.Range("E22").Value = "  +1.51%  "
$ws.Range("E23").Value = "  +1.15%  "
$ws.Range("D24").Value = "'13.61"
$ws.Range("E24").Value = "  -1.16%  "
$ws.Range("D25").Value = "'84.97"
$ws.Range("E25").Value = "  +1.93%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Value = "'2.78"
$ws.Range("E27").Value = "  -2.99%  "
$ws.Range("E28").Value = "  +1.34%  "
$ws.Range("D29").Value = "'7.36"
$ws.Range("E29").Value = "  +7.74%  "
$ws.Range("E30").Value = "  +3.25%  "
$ws.Range("E31").Value = "  -4.90%  "
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").Value = "'26.70"
$ws.Range("E33").Value = "  +1.85%  "
$ws.Range("D34").Value = "'2.64"
$ws.Range("E34").Value = "  -4.63%  "
$ws.Range("E35").Value = "  +0.64%  "
$ws.Range("D36").Value = "'5.94"
$ws.Range("E36").Value = "  -0.79%  "
$ws.Range("D37").Value = "'52.43"
$ws.Range("E37").Value = "  -3.29%  "
$ws.Range("D38").Value = "0.0₃0743"
$ws.Range("E38").Value = "  +4.08%  "
$ws.Range("D39").Value = "'453.08"
$ws.Range("E39").Value = "  +1.63%  "
$ws.Range("E40").Value = "  +3.69%  "
$ws.Range("D41").Value = "'0.0395"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E42").Value = "  -0.59%  "
$ws.Range("D43").Value = "'8.33"
$ws.Range("E43").Value = "  -1.35%  "
$ws.Range("D44").Value = "2.857.28"
$ws.Range("E44").Value = "  +0.87%  "
$ws.Range("E45").Value = "  -0.73%  "
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("E47").Value = "  +4.50%  "
$ws.Range("D48").Value = "'26.42"
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("D51").Value = "'120.14"
$ws.Range("E51").Value = "  +1.84%  "
